# Update "想去人数" (F column) counts on sheet "展览" (sheet1)
$ws1 = $excel.ActiveWorkbook.Worksheets.Item("展览")
$ws1.Range("F4").Value = 586
$ws1.Range("F6").Value = 29
$ws1.Range("F7").Value = 1926
$ws1.Range("F8").Value = 5281
$ws1.Range("F9").Value = 1472
$ws1.Range("F11").Value = 3071
$ws1.Range("F14").Value = 1279
$ws1.Range("F15").Value = 4206
$ws1.Range("F16").Value = 1004
$ws1.Range("F17").Value = 886
$ws1.Range("F18").Value = 1650
$ws1.Range("F19").Value = 2594
$ws1.Range("F22").Value = 128
$ws1.Range("F23").Value = 143
$ws1.Range("F24").Value = 960
$ws1.Range("F25").Value = 287
$ws1.Range("F29").Value = 1069
$ws1.Range("F30").Value = 357
$ws1.Range("F32").Value = 119
$ws1.Range("F34").Value = 237
$ws1.Range("F35").Value = 1636
$ws1.Range("F36").Value = 2157
$ws1.Range("F37").Value = 1006
$ws1.Range("F40").Value = 598
$ws1.Range("F41").Value = 270
$ws1.Range("F45").Value = 308
$ws1.Range("F46").Value = 205
$ws1.Range("F47").Value = 130

# Update "想去人数" (F column) counts on sheet "本地生活" (sheet3)
$ws3 = $excel.ActiveWorkbook.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 727

# Update "想去人数" (F column) counts on sheet "全部类型" (sheet4)
$ws4 = $excel.ActiveWorkbook.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 727
$ws4.Range("F6").Value = 586
$ws4.Range("F7").Value = 29
$ws4.Range("F8").Value = 1926
$ws4.Range("F9").Value = 5281
$ws4.Range("F10").Value = 1472
$ws4.Range("F13").Value = 3071
$ws4.Range("F15").Value = 1279
$ws4.Range("F16").Value = 4206
$ws4.Range("F17").Value = 1004
$ws4.Range("F18").Value = 1650
$ws4.Range("F20").Value = 2594
$ws4.Range("F25").Value = 143
$ws4.Range("F27").Value = 960
$ws4.Range("F28").Value = 287
$ws4.Range("F32").Value = 1069
$ws4.Range("F33").Value = 357
$ws4.Range("F36").Value = 1636
$ws4.Range("F37").Value = 2157
$ws4.Range("F38").Value = 1006
$ws4.Range("F42").Value = 598
$ws4.Range("F43").Value = 270
$ws4.Range("F46").Value = 308
$ws4.Range("F47").Value = 205
$ws4.Range("F48").Value = 130
